# Edit script: add two new rows (13 and 14) of translation-time data to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 13: gakumas lilja white night white wish / gpt-4o
$ws.Range("A13").Value = "gakumas lilja white night white wish"
$ws.Range("B13").Value = "gpt-4o"
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 9
$ws.Range("E13").Value = 18
$ws.Range("F13").Value = 242
$ws.Range("G13").Value = 5706
$ws.Range("H13").Formula = "=G13/F13"
$ws.Range("H13").NumberFormat = "0.00"
$ws.Range("I13").Value = 3
$ws.Range("J13").Value = 0.42
$ws.Range("K13").Formula = "=J13/G13"
$ws.Range("K13").NumberFormat = "0.000000"
$ws.Range("L13").Formula = "=J13/F13"

# Row 14: gakumas / gpt-4.1-mini
$ws.Range("A14").Value = "gakumas"
$ws.Range("B14").Value = "gpt-4.1-mini"
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 53
$ws.Range("E14").Value = 38
$ws.Range("F14").Value = 690
$ws.Range("G14").Value = 16039
$ws.Range("H14").Formula = "=G14/F14"
$ws.Range("H14").NumberFormat = "0.00"
$ws.Range("I14").Value = 3
$ws.Range("J14").Formula = "=1.49-0.25"
$ws.Range("K14").Formula = "=J14/G14"
$ws.Range("K14").NumberFormat = "0.000000"
$ws.Range("L14").Formula = "=J14/F14"

# Update selection to match final state
$ws.Range("R16").Select()

$wb.Save()
